$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlink currently on C2 before changing the cell value
$ws.Range("C2").Hyperlinks.Delete()

# Update C2's value (the picture link for Leopard Gecko) to the new URL
$ws.Range("C2").Value = "https://geckoadvice.com/wp-content/uploads/2022/03/Leopard-Gecko-Climbing.jpg"

# Delete column D (the Alignment column) entirely
$ws.Range("D1:D5").EntireColumn.Delete()

# Update the active cell selection to D1 (matches target workbook state)
$ws.Range("D1").Select()
